# This script applies a "shift down by one row + append new last row" update
# to the weekly "Bruselas (repollito)" price sheet, and inserts one brand-new
# data row (the latest week) at the top of the data block (row 4).
#
# Pattern (confirmed against the target diff):
#   - Rows 2 and 3 are untouched.
#   - Row 4 becomes a brand-new record (new date/volume/prices).
#   - Every row from 5 to 63 takes on the values that the row immediately
#     above it held *before* this edit (i.e. the whole block from old row 4
#     through old row 62 is shifted down by one row).
#   - A brand-new row 64 is appended, holding the values that used to be in
#     old row 63 (so nothing from the original data set is lost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Snapshot the current (pre-edit) contents of rows 4 through 63 (whole
#    row, columns A:R) before anything is overwritten.
$snapshot = @()
for ($r = 4; $r -le 63; $r++) {
    $snapshot += ,$ws.Range("A$r`:R$r").Value2
}

# 2. Push that snapshot down by one row: old row N (N = 4..63) becomes the
#    new row N+1 (new rows 5..64).
for ($r = 5; $r -le 64; $r++) {
    $idx = $r - 5
    $ws.Range("A$r`:R$r").Value2 = $snapshot[$idx]
}

# 3. Make sure the newly created row 64 has the same date number format as
#    the rest of the "Fecha" column (column D).
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(63, 4).NumberFormat

# 4. Overwrite row 4 with the brand-new record.
$ws.Cells.Item(4, 4).Value = 45043
$ws.Cells.Item(4, 10).Value = 220
$ws.Cells.Item(4, 11).Value = 18000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 18909
$ws.Cells.Item(4, 16).Value = 1261

Write-Host "Applied weekly shift update; dimension now A1:R64"
